# Updated cryptos list values (price + 1h volume change), matching the
# upstream scraper commit. Rows 27/28 and 42/43 also swap rank order
# (coin name + link), so B/C for those rows are rewritten too.
$updates = @(
    @{ Cell = "D2"; Value = '63.361.39'; ForceText = 0 }
    @{ Cell = "E2"; Value = '  -1.09%  '; ForceText = 0 }
    @{ Cell = "D3"; Value = '3.257.18'; ForceText = 0 }
    @{ Cell = "E3"; Value = '  +3.29%  '; ForceText = 0 }
    @{ Cell = "E4"; Value = '  -0.14%  '; ForceText = 0 }
    @{ Cell = "D5"; Value = '594.20'; ForceText = 1 }
    @{ Cell = "E5"; Value = '  -1.60%  '; ForceText = 0 }
    @{ Cell = "D6"; Value = '141.26'; ForceText = 1 }
    @{ Cell = "E6"; Value = '  -1.73%  '; ForceText = 0 }
    @{ Cell = "E7"; Value = '  -0.06%  '; ForceText = 0 }
    @{ Cell = "D8"; Value = '3.255.65'; ForceText = 0 }
    @{ Cell = "E8"; Value = '  +3.41%  '; ForceText = 0 }
    @{ Cell = "E9"; Value = '  -0.80%  '; ForceText = 0 }
    @{ Cell = "E10"; Value = '  -1.15%  '; ForceText = 0 }
    @{ Cell = "D11"; Value = '5.34'; ForceText = 1 }
    @{ Cell = "E11"; Value = '  -0.84%  '; ForceText = 0 }
    @{ Cell = "D12"; Value = '0.466'; ForceText = 1 }
    @{ Cell = "E12"; Value = '  -0.21%  '; ForceText = 0 }
    @{ Cell = "E13"; Value = '  -2.39%  '; ForceText = 0 }
    @{ Cell = "D14"; Value = '34.68'; ForceText = 1 }
    @{ Cell = "E14"; Value = '  -1.00%  '; ForceText = 0 }
    @{ Cell = "D15"; Value = '3.789.79'; ForceText = 0 }
    @{ Cell = "E15"; Value = '  +3.11%  '; ForceText = 0 }
    @{ Cell = "E16"; Value = '  +0.05%  '; ForceText = 0 }
    @{ Cell = "D17"; Value = '3.250.84'; ForceText = 0 }
    @{ Cell = "E17"; Value = '  +3.03%  '; ForceText = 0 }
    @{ Cell = "D18"; Value = '63.397.86'; ForceText = 0 }
    @{ Cell = "E18"; Value = '  -1.15%  '; ForceText = 0 }
    @{ Cell = "D19"; Value = '6.79'; ForceText = 1 }
    @{ Cell = "E19"; Value = '  -0.90%  '; ForceText = 0 }
    @{ Cell = "D20"; Value = '477.19'; ForceText = 1 }
    @{ Cell = "E20"; Value = '  -3.00%  '; ForceText = 0 }
    @{ Cell = "D21"; Value = '14.22'; ForceText = 1 }
    @{ Cell = "E21"; Value = '  -3.43%  '; ForceText = 0 }
    @{ Cell = "D22"; Value = '0.728'; ForceText = 1 }
    @{ Cell = "E22"; Value = '  +2.28%  '; ForceText = 0 }
    @{ Cell = "D23"; Value = '7.98'; ForceText = 1 }
    @{ Cell = "E23"; Value = '  +4.24%  '; ForceText = 0 }
    @{ Cell = "D24"; Value = '83.99'; ForceText = 1 }
    @{ Cell = "E24"; Value = '  -4.41%  '; ForceText = 0 }
    @{ Cell = "E25"; Value = '  -0.24%  '; ForceText = 0 }
    @{ Cell = "B27"; Value = 'NEARProtocol'; ForceText = 0 }
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; ForceText = 0 }
    @{ Cell = "D27"; Value = '7.50'; ForceText = 1 }
    @{ Cell = "E27"; Value = '  +7.17%  '; ForceText = 0 }
    @{ Cell = "B28"; Value = 'PancakeSwap'; ForceText = 0 }
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; ForceText = 0 }
    @{ Cell = "D28"; Value = '2.74'; ForceText = 1 }
    @{ Cell = "E28"; Value = '  -1.03%  '; ForceText = 0 }
    @{ Cell = "D29"; Value = '8.09'; ForceText = 1 }
    @{ Cell = "E29"; Value = '  -1.54%  '; ForceText = 0 }
    @{ Cell = "D30"; Value = '2.13'; ForceText = 1 }
    @{ Cell = "E30"; Value = '  +3.31%  '; ForceText = 0 }
    @{ Cell = "D31"; Value = '27.73'; ForceText = 1 }
    @{ Cell = "E31"; Value = '  +0.48%  '; ForceText = 0 }
    @{ Cell = "E33"; Value = '  -2.57%  '; ForceText = 0 }
    @{ Cell = "D34"; Value = '2.57'; ForceText = 1 }
    @{ Cell = "E34"; Value = '  -3.55%  '; ForceText = 0 }
    @{ Cell = "D35"; Value = '1.10'; ForceText = 1 }
    @{ Cell = "E35"; Value = '  -0.99%  '; ForceText = 0 }
    @{ Cell = "D36"; Value = '5.92'; ForceText = 1 }
    @{ Cell = "E36"; Value = '  -1.73%  '; ForceText = 0 }
    @{ Cell = "D37"; Value = '52.97'; ForceText = 1 }
    @{ Cell = "E37"; Value = '  +0.43%  '; ForceText = 0 }
    @{ Cell = "D38"; Value = '0.0₃0720'; ForceText = 0 }
    @{ Cell = "E38"; Value = '  -3.37%  '; ForceText = 0 }
    @{ Cell = "D39"; Value = '0.0394'; ForceText = 1 }
    @{ Cell = "E39"; Value = '  -1.04%  '; ForceText = 0 }
    @{ Cell = "D40"; Value = '422.37'; ForceText = 1 }
    @{ Cell = "E40"; Value = '  -3.04%  '; ForceText = 0 }
    @{ Cell = "D41"; Value = '3.002.45'; ForceText = 0 }
    @{ Cell = "E41"; Value = '  +2.02%  '; ForceText = 0 }
    @{ Cell = "B42"; Value = 'dogwifhat'; ForceText = 0 }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; ForceText = 0 }
    @{ Cell = "D42"; Value = '2.77'; ForceText = 1 }
    @{ Cell = "E42"; Value = '  -6.59%  '; ForceText = 0 }
    @{ Cell = "B43"; Value = 'Cosmos'; ForceText = 0 }
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText = 0 }
    @{ Cell = "D43"; Value = '8.41'; ForceText = 1 }
    @{ Cell = "E43"; Value = '  +1.22%  '; ForceText = 0 }
    @{ Cell = "D44"; Value = '0.112'; ForceText = 1 }
    @{ Cell = "E44"; Value = '  -6.30%  '; ForceText = 0 }
    @{ Cell = "D45"; Value = '0.270'; ForceText = 1 }
    @{ Cell = "E45"; Value = '  +4.13%  '; ForceText = 0 }
    @{ Cell = "D46"; Value = '2.17'; ForceText = 1 }
    @{ Cell = "E46"; Value = '  -1.18%  '; ForceText = 0 }
    @{ Cell = "D48"; Value = '25.99'; ForceText = 1 }
    @{ Cell = "E48"; Value = '  +0.09%  '; ForceText = 0 }
    @{ Cell = "D49"; Value = '2.33'; ForceText = 1 }
    @{ Cell = "E49"; Value = '  -3.12%  '; ForceText = 0 }
    @{ Cell = "D50"; Value = '0.115'; ForceText = 1 }
    @{ Cell = "E50"; Value = '  +0.28%  '; ForceText = 0 }
    @{ Cell = "D51"; Value = '33.97'; ForceText = 1 }
    @{ Cell = "E51"; Value = '  +9.37%  '; ForceText = 0 }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText -eq 1) {
        # Force text so Excel doesn't re-interpret a numeric-looking price
        # (e.g. "594.20") as a number and drop the trailing zero / change type.
        $range.Value = "'" + $u.Value
    } else {
        $range.Value = $u.Value
    }
}

Write-Host "Applied $($updates.Count) cell updates"
